$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be read/stored as plain text so values such as
# "1.002", "327.20" or "92.00" keep their exact original formatting instead of
# being auto-converted into numbers (which would drop trailing/insignificant
# zeros, e.g. "92.00" -> 92). Column E (percentages, padded with spaces) and
# column B/C (names/links) are never mistaken for numbers, so they need no
# special handling.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.841.22'
$ws.Range("E2").Value = '  +1.01%  '

$ws.Range("D3").Value = '1.755.11'
$ws.Range("E3").Value = '  +0.00%  '

$ws.Range("D4").Value = '1.002'

$ws.Range("D5").Value = '327.20'
$ws.Range("E5").Value = '  +0.69%  '

$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").Value = '0.4629'
$ws.Range("E7").Value = '  +1.18%  '

$ws.Range("D8").Value = '0.3487'
$ws.Range("E8").Value = '  -2.06%  '

$ws.Range("D9").Value = '41.98'
$ws.Range("E9").Value = '  +1.14%  '

$ws.Range("D10").Value = '0.07348'
$ws.Range("E10").Value = '  -1.41%  '

$ws.Range("D11").Value = '1.078'
$ws.Range("E11").Value = '  -0.66%  '

$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.02%  '

$ws.Range("D13").Value = '20.46'
$ws.Range("E13").Value = '  -1.57%  '

$ws.Range("D14").Value = '5.972'
$ws.Range("E14").Value = '  -0.59%  '

$ws.Range("D15").Value = '7.134'
$ws.Range("E15").Value = '  -0.40%  '

$ws.Range("D16").Value = '1.756.57'
$ws.Range("E16").Value = '  +0.12%  '

$ws.Range("D17").Value = '92.00'
$ws.Range("E17").Value = '  -1.61%  '

$ws.Range("E18").Value = '  -0.01%  '

$ws.Range("D19").Value = '0.06405'
$ws.Range("E19").Value = '  -0.11%  '

$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.03%  '

$ws.Range("D21").Value = '16.77'
$ws.Range("E21").Value = '  -1.68%  '

$ws.Range("D22").Value = '5.748'
$ws.Range("E22").Value = '  +0.14%  '

$ws.Range("D23").Value = '27.854.97'
$ws.Range("E23").Value = '  +0.90%  '

$ws.Range("D24").Value = '11.12'
$ws.Range("E24").Value = '  -0.90%  '

$ws.Range("D25").Value = '2.152'
$ws.Range("E25").Value = '  +3.95%  '

$ws.Range("D26").Value = '161.59'
$ws.Range("E26").Value = '  -2.22%  '

$ws.Range("D27").Value = '19.98'
$ws.Range("E27").Value = '  -0.70%  '

$ws.Range("D28").Value = '1.956.68'
$ws.Range("E28").Value = '  +0.08%  '

$ws.Range("D29").Value = '2.146'
$ws.Range("E29").Value = '  +0.47%  '

$ws.Range("D30").Value = '122.60'
$ws.Range("E30").Value = '  -2.28%  '

$ws.Range("D31").Value = '1.066'
$ws.Range("E31").Value = '  -1.00%  '

$ws.Range("D32").Value = '0.09294'
$ws.Range("E32").Value = '  +0.79%  '

$ws.Range("D33").Value = '3.654'
$ws.Range("E33").Value = '  -0.27%  '

$ws.Range("D34").Value = '5.529'
$ws.Range("E34").Value = '  +0.36%  '

$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '11.60'
$ws.Range("E35").Value = '  -1.19%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.02258'
$ws.Range("E36").Value = '  -0.68%  '

$ws.Range("D37").Value = '0.06056'
$ws.Range("E37").Value = '  +0.62%  '

$ws.Range("D38").Value = '0.2059'
$ws.Range("E38").Value = '  -1.15%  '

$ws.Range("D39").Value = '4.885'
$ws.Range("E39").Value = '  -0.79%  '

$ws.Range("D40").Value = '0.6149'
$ws.Range("E40").Value = '  -1.83%  '

$ws.Range("D41").Value = '1.176'
$ws.Range("E41").Value = '  -0.47%  '

$ws.Range("D42").Value = '7.736'
$ws.Range("E42").Value = '  -0.27%  '

$ws.Range("D43").Value = '1.345'
$ws.Range("E43").Value = '  -2.74%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '12.99'
$ws.Range("E44").Value = '  -1.80%  '

$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").Value = '3.730'
$ws.Range("E45").Value = '  +0.35%  '

$ws.Range("E46").Value = '  -1.59%  '

$ws.Range("D47").Value = '122.54'
$ws.Range("E47").Value = '  +0.54%  '

$ws.Range("D48").Value = '1.920'
$ws.Range("E48").Value = '  -0.77%  '

$ws.Range("E49").Value = '  -1.53%  '

$ws.Range("E50").Value = '  -1.12%  '

$ws.Range("D51").Value = '71.96'
$ws.Range("E51").Value = '  -0.13%  '
